# Add two generic inputs to wally (two GPIO pins already on board,
# now wiring up a generic 3-pin header part - UJC-HP2-3-SMT-TR - instead of
# the old USB-C connector 217175-0001, and marking several "Other" parts
# as "Added to order".

$wb = $excel.ActiveWorkbook
$wsOther = $wb.Worksheets.Item("Other")
$wsJlc = $wb.Worksheets.Item("JLC")

# --- Update the USB/connector part on row 5 of "Other" to the new part ---
$wsOther.Range("C5").Value2 = "UJC-HP2-3-SMT-TR"
$wsOther.Range("D5").Value2 = "https://www.mouser.co.uk/ProductDetail/CUI-Devices/UJC-HP2-3-SMT-TR?qs=HoCaDK9Nz5cglCCyoWNzZg%3D%3D"
$wsOther.Range("F5").Value2 = 0.464

# --- Add a new "Added to order" column to the Table13 list object ---
$lo = $wsOther.ListObjects.Item("Table13")
$newCol = $lo.ListColumns.Add()
$wsOther.Range("G1").Value2 = "Added to order"

# --- Mark rows 2 through 13 as added to order ---
for ($r = 2; $r -le 13; $r++) {
    $wsOther.Cells.Item($r, 7).Value2 = "y"
}

# --- Switch active sheet / selection to mirror the author's last view ---
$wsOther.Activate()
$wsOther.Range("I19").Select()
